$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is entered first so its text becomes the next shared-string entry,
# followed by row 4's text, followed by row 6's text (matches target order).
$ws.Range("A5").Value = "VETERINARY SYSTEM FEATURES AND BASIS"
$ws.Range("B5").Value = 43726
$ws.Range("B5").NumberFormat = "d-mmm"
$ws.Range("C5").Value = 43727
$ws.Range("C5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = 0.91666666666666663
$ws.Range("D5").NumberFormat = "h:mm AM/PM"
$ws.Range("E5").Value = 0.083333333333333329
$ws.Range("E5").NumberFormat = "h:mm AM/PM"

$ws.Range("A4").Value = "CLIENT FINDING"
$ws.Range("B4").Value = 43731
$ws.Range("B4").NumberFormat = "d-mmm"
$ws.Range("C4").Value = 43731
$ws.Range("C4").NumberFormat = "d-mmm"
$ws.Range("D4").Value = 0.66666666666666663
$ws.Range("D4").NumberFormat = "h:mm AM/PM"
$ws.Range("E4").Value = 0.75
$ws.Range("E4").NumberFormat = "h:mm AM/PM"

$ws.Range("A6").Value = "NEW ADDITIONAL FEATURES RESEARCH"
$ws.Range("B6").Value = 43733
$ws.Range("B6").NumberFormat = "d-mmm"
$ws.Range("C6").Value = 43734
$ws.Range("C6").NumberFormat = "d-mmm"
$ws.Range("D6").Value = 0.89583333333333337
$ws.Range("D6").NumberFormat = "h:mm AM/PM"
$ws.Range("E6").Value = 0.95208333333333339
$ws.Range("E6").NumberFormat = "h:mm AM/PM"

$ws.Range("E17").Select()
